$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that may contain HYPERLINK(...) formulas that need a friendly
# second argument added (the "Beteckning" value from column A).
$linkCols = @(19, 20, 21, 22, 23, 24, 25)  # S, T, U, V, W, X, Y

$firstRow = 2
$lastRow = 250

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # Update the "Förändrad" (changed) date in column C (3) to the new value.
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value = 45186
    }

    # Grab the "Beteckning" value in column A to use as the friendly
    # hyperlink display text.
    $name = $ws.Cells.Item($r, 1).Value2

    foreach ($c in $linkCols) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -like '*HYPERLINK(*' -and $f -notlike '*,*') {
                $newF = $f -replace '\)$', (', "' + $name + '")')
                $cell.Formula = $newF
            }
        }
    }
}
